# Regenerate save_data to use K (column G) instead of Strike# values.
# This writes the newly calculated "K" values (s_vals) into column G
# for each data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K value (column G), per the recalculated
# std/mean based s_vals.
$kValues = @{
    2  = 0
    4  = 0
    5  = 1
    6  = 1
    7  = 2
    8  = 1
    9  = 1
    10 = 1
    11 = 0
    12 = 1
    13 = 0
    14 = 0
    15 = 3
    16 = 0
    17 = 1
    18 = 2
    19 = 1
    20 = 3
    21 = 1
    22 = 2
    23 = 1
    24 = 2
    25 = 1
    26 = 0
    27 = 0
    28 = 2
    29 = 1
    30 = 2
    31 = 1
    32 = 2
    33 = 2
    34 = 1
    35 = 1
    36 = 1
    37 = 4
    38 = 0
    39 = 2
    40 = 1
    41 = 2
    42 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
